$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 1983.0834
$ws.Range("I62").Value = 2007.5555
$ws.Range("K62").Value = 2007.5555
$ws.Range("M62").Value = -1383.5555
$ws.Range("H65").Value = 1983.0834
$ws.Range("I65").Value = 2007.5555
$ws.Range("K65").Value = 10037.7775
$ws.Range("M65").Value = -6917.7775
$ws.Range("H74").Value = 5455.4443
$ws.Range("I74").Value = 3900
$ws.Range("K74").Value = 3900
$ws.Range("M74").Value = -2964
$ws.Range("H77").Value = 5455.4443
$ws.Range("I77").Value = 3900
$ws.Range("K77").Value = 19500
$ws.Range("M77").Value = -14820
$ws.Range("H106").Value = 3109.8
$ws.Range("I106").Value = 3109.8
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 3109.8
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -2478.8
$ws.Range("N106").ClearContents()
$ws.Range("H113").Value = 68593.664
$ws.Range("I113").Value = 201841
$ws.Range("J113").Value = 1970
$ws.Range("K113").Value = 201841
$ws.Range("L113").Value = 1970
$ws.Range("M113").Value = -198587
$ws.Range("N113").Value = -8478
$ws.Range("H129").Value = 2585.1
$ws.Range("J129").Value = 966.39624
$ws.Range("L129").Value = 2899.18872
$ws.Range("N129").Value = -12899.18872
$ws.Range("H132").Value = 7149104.5
$ws.Range("I132").Value = 7582303.5
$ws.Range("K132").Value = 22746910.5
$ws.Range("M132").Value = -22744380.5
$ws.Range("H137").Value = 1498.05
$ws.Range("I137").Value = 1438.3
$ws.Range("K137").Value = 4314.9
$ws.Range("M137").Value = -1764.9
$ws.Range("H141").Value = 2289.087
$ws.Range("I141").Value = 1364
$ws.Range("K141").Value = 4092
$ws.Range("M141").Value = 1088

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 27184.643
$ws.Range("I32").Value = 4569.8726
$ws.Range("K32").Value = 4569.8726
$ws.Range("M32").Value = -4282.8726
$ws.Range("H45").Value = 1886.2609
$ws.Range("I45").Value = 1727.1333
$ws.Range("J45").Value = 2184.625
$ws.Range("K45").Value = 1727.1333
$ws.Range("L45").Value = 2184.625
$ws.Range("M45").Value = -1350.1333
$ws.Range("N45").Value = -2938.625
$ws.Range("H69").Value = 41026.375
$ws.Range("J69").Value = 41026.375
$ws.Range("L69").Value = 41026.375
$ws.Range("N69").Value = -42524.375
$ws.Range("H72").Value = 41026.375
$ws.Range("J72").Value = 41026.375
$ws.Range("L72").Value = 123079.125
$ws.Range("N72").Value = -130567.125
$ws.Range("H74").Value = 1856.3846
$ws.Range("I74").Value = 1050.2632
$ws.Range("J74").Value = 2622.2
$ws.Range("K74").Value = 1050.2632
$ws.Range("L74").Value = 2622.2
$ws.Range("M74").Value = -176.2632000000001
$ws.Range("N74").Value = -4370.2
$ws.Range("H77").Value = 1856.3846
$ws.Range("I77").Value = 1050.2632
$ws.Range("J77").Value = 2622.2
$ws.Range("K77").Value = 5251.316000000001
$ws.Range("L77").Value = 13111
$ws.Range("M77").Value = -883.3160000000007
$ws.Range("N77").Value = -21847
$ws.Range("H122").Value = 2948.8635
$ws.Range("I122").Value = 3367.8125
$ws.Range("J122").Value = 1831.6666
$ws.Range("K122").Value = 10103.4375
$ws.Range("L122").Value = 5494.9998
$ws.Range("M122").Value = -7653.4375
$ws.Range("N122").Value = -10394.9998

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 166991.67
$ws.Range("I94").Value = 200290
$ws.Range("K94").Value = 200290
$ws.Range("M94").Value = -199839
$ws.Range("H105").Value = 78857.30499999999
$ws.Range("I105").Value = 51813.95
$ws.Range("J105").Value = 169001.83
$ws.Range("K105").Value = 51813.95
$ws.Range("L105").Value = 169001.83
$ws.Range("M105").Value = -50066.95
$ws.Range("N105").Value = -172495.83
$ws.Range("H134").Value = 2765.8235
$ws.Range("I134").Value = 2794.2144
$ws.Range("J134").Value = 2633.3333
$ws.Range("K134").Value = 8382.643199999999
$ws.Range("L134").Value = 7899.999899999999
$ws.Range("M134").Value = -5847.643199999999
$ws.Range("N134").Value = -12969.9999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 17472.828
$ws.Range("I31").Value = 41114.56
$ws.Range("J31").Value = 2317.8718
$ws.Range("K31").Value = 41114.56
$ws.Range("L31").Value = 2317.8718
$ws.Range("M31").Value = -40819.56
$ws.Range("N31").Value = -2907.8718
$ws.Range("H34").Value = 17472.828
$ws.Range("I34").Value = 41114.56
$ws.Range("J34").Value = 2317.8718
$ws.Range("K34").Value = 41114.56
$ws.Range("L34").Value = 2317.8718
$ws.Range("M34").Value = -40912.56
$ws.Range("N34").Value = -2721.8718
$ws.Range("H58").Value = 12815.608
$ws.Range("I58").Value = 2079.3845
$ws.Range("J58").Value = 26772.7
$ws.Range("K58").Value = 2079.3845
$ws.Range("L58").Value = 26772.7
$ws.Range("M58").Value = -1876.3845
$ws.Range("N58").Value = -27178.7
$ws.Range("H107").Value = 4883.875
$ws.Range("I107").Value = 9870.909
$ws.Range("J107").Value = 664.0769
$ws.Range("K107").Value = 9870.909
$ws.Range("L107").Value = 664.0769
$ws.Range("M107").Value = -7950.909
$ws.Range("N107").Value = -4504.0769
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H132").Value = 3700.611
$ws.Range("I132").Value = 4273.727
$ws.Range("K132").Value = 12821.181
$ws.Range("M132").Value = -10291.181
$ws.Range("H134").Value = 1556.1305
$ws.Range("I134").Value = 1315.3889
$ws.Range("K134").Value = 3946.1667
$ws.Range("M134").Value = -1411.1667
$ws.Range("H136").Value = 12815.608
$ws.Range("I136").Value = 2079.3845
$ws.Range("J136").Value = 26772.7
$ws.Range("K136").Value = 6238.1535
$ws.Range("L136").Value = 80318.10000000001
$ws.Range("M136").Value = -3688.1535
$ws.Range("N136").Value = -85418.10000000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 557.4
$ws.Range("I23").Value = 115
$ws.Range("J23").Value = 625.46155
$ws.Range("K23").Value = 345
$ws.Range("L23").Value = 1876.38465
$ws.Range("M23").Value = -110
$ws.Range("N23").Value = -2346.38465
$ws.Range("H68").Value = 1812.84
$ws.Range("I68").Value = 1245.1034
$ws.Range("J68").Value = 2170.761
$ws.Range("K68").Value = 3735.3102
$ws.Range("L68").Value = 6512.282999999999
$ws.Range("M68").Value = -2924.3102
$ws.Range("N68").Value = -8134.282999999999
$ws.Range("H71").Value = 1812.84
$ws.Range("I71").Value = 1245.1034
$ws.Range("J71").Value = 2170.761
$ws.Range("K71").Value = 11205.9306
$ws.Range("L71").Value = 19536.849
$ws.Range("M71").Value = -7149.9306
$ws.Range("N71").Value = -27648.849
$ws.Range("H131").Value = 1218.7273
$ws.Range("J131").Value = 1187.5952
$ws.Range("L131").Value = 3562.7856
$ws.Range("N131").Value = -13642.7856

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 52633824
$ws.Range("I97").Value = 71430990
$ws.Range("K97").Value = 71430990
$ws.Range("M97").Value = -71430494
$ws.Range("H102").Value = 355937.88
$ws.Range("I102").Value = 3965.1
$ws.Range("K102").Value = 3965.1
$ws.Range("M102").Value = -2343.1
$ws.Range("H122").Value = 3506.3
$ws.Range("I122").Value = 4073.0908
$ws.Range("J122").Value = 2813.5557
$ws.Range("K122").Value = 12219.2724
$ws.Range("L122").Value = 8440.667099999999
$ws.Range("M122").Value = -9769.2724
$ws.Range("N122").Value = -13340.6671
$ws.Range("H123").Value = 9326
$ws.Range("J123").Value = 9326
$ws.Range("L123").Value = 9326
$ws.Range("N123").Value = -14226
$ws.Range("H126").Value = 5884812
$ws.Range("I126").Value = 3396
$ws.Range("J126").Value = 11766228
$ws.Range("K126").Value = 10188
$ws.Range("L126").Value = 35298684
$ws.Range("M126").Value = -7718
$ws.Range("N126").Value = -35303624
$ws.Range("H132").Value = 3463.9
$ws.Range("I132").Value = 3182.2222
$ws.Range("K132").Value = 9546.6666
$ws.Range("M132").Value = -7016.6666

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2500
$ws.Range("I61").Value = 2500
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2500
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2298
$ws.Range("N61").ClearContents()
$ws.Range("H93").Value = 1413.1666
$ws.Range("I93").Value = 1483.7693
$ws.Range("J93").Value = 1229.6
$ws.Range("K93").Value = 1483.7693
$ws.Range("L93").Value = 1229.6
$ws.Range("M93").Value = -235.7692999999999
$ws.Range("N93").Value = -3725.6
$ws.Range("H113").Value = 2500
$ws.Range("I113").Value = 2500
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2500
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -330
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 3431.8333
$ws.Range("I122").Value = 3431.8333
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 10295.4999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -7845.499899999999
$ws.Range("N122").ClearContents()
$ws.Range("H136").Value = 2485.6428
$ws.Range("I136").Value = 1889.4
$ws.Range("J136").Value = 3976.25
$ws.Range("K136").Value = 5668.200000000001
$ws.Range("L136").Value = 11928.75
$ws.Range("M136").Value = -3118.200000000001
$ws.Range("N136").Value = -17028.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("H107").Value = 64456.156
$ws.Range("I107").Value = 22825.223
$ws.Range("J107").Value = 101924
$ws.Range("K107").Value = 68475.66900000001
$ws.Range("L107").Value = 305772
$ws.Range("M107").Value = -66555.66900000001
$ws.Range("N107").Value = -309612
$ws.Range("H132").Value = 6202.125
$ws.Range("I132").Value = 6381.0713
$ws.Range("J132").Value = 4949.5
$ws.Range("K132").Value = 19143.2139
$ws.Range("L132").Value = 14848.5
$ws.Range("M132").Value = -16613.2139
$ws.Range("N132").Value = -19908.5
$ws.Range("H136").Value = 1365.8096
$ws.Range("I136").Value = 960.53845
$ws.Range("J136").Value = 2024.375
$ws.Range("K136").Value = 2881.61535
$ws.Range("L136").Value = 6073.125
$ws.Range("M136").Value = -331.61535
$ws.Range("N136").Value = -11173.125
